# 2.1.1.1e.xlsx — add a "2020" data column (N) and refresh several existing
# data points for 2018/2019, matching the upstream OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New column N: header (year 2020) + data for rows 5-14, plus the
#    thin bottom-border spacer cell in row 3 (mirrors column M's shape).
# ---------------------------------------------------------------------

# Row 3 spacer cell (style only, no value) — copy format from the existing
# bottom-border row-3 cells (e.g. A3, which carries the same border/font).
$ws.Range("A14").Copy()
$ws.Range("N3").PasteSpecial(-4122)

# Row 4 header "2020" — copy the year-header format from D4 (style 5).
$ws.Range("D4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 2020

# Rows 5-13 — numeric data cells, same number format/border as column D
# (style 10).
$ws.Range("D5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 68.5

$ws.Range("D6").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N6").Value = 106.7

$ws.Range("D7").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$ws.Range("N7").Value = 53.2

$ws.Range("D8").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("N8").Value = 49.6

$ws.Range("D9").Copy()
$ws.Range("N9").PasteSpecial(-4122)
$ws.Range("N9").Value = 108.9

$ws.Range("D10").Copy()
$ws.Range("N10").PasteSpecial(-4122)
$ws.Range("N10").Value = 107.8

$ws.Range("D11").Copy()
$ws.Range("N11").PasteSpecial(-4122)
$ws.Range("N11").Value = 155.7

$ws.Range("D12").Copy()
$ws.Range("N12").PasteSpecial(-4122)
$ws.Range("N12").Value = 25.9

$ws.Range("D13").Copy()
$ws.Range("N13").PasteSpecial(-4122)
$ws.Range("N13").Value = 103.5

# Row 14 — last row uses the thicker/bottom-border style (style 12),
# mirrored from M14.
$ws.Range("M14").Copy()
$ws.Range("N14").PasteSpecial(-4122)
$ws.Range("N14").Value = 11

# ---------------------------------------------------------------------
# 2. Refresh existing 2018 (L) / 2019 (M) figures that moved with the
#    new data release.
# ---------------------------------------------------------------------

$ws.Range("M5").Value = 68.400000000000006
$ws.Range("M6").Value = 108.2
$ws.Range("M7").Value = 51.7
$ws.Range("M8").Value = 97.7

$ws.Range("L9").Value = 105.6
$ws.Range("M9").Value = 106.7

$ws.Range("M10").Value = 124.2
$ws.Range("M11").Value = 138.80000000000001

$ws.Range("L12").Value = 27.1
$ws.Range("M12").Value = 33.9

$ws.Range("M13").Value = 96
$ws.Range("M14").Value = 7.7

# ---------------------------------------------------------------------
# 3. Page setup — paper size / orientation were set via Page Setup.
# ---------------------------------------------------------------------

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
